$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the state name entries (first letter uppercase)
$ws.Range("G1").Value = "New Mexico"
$ws.Range("G3").Value = "Uttah"
$ws.Range("G4").Value = "Alabama"
$ws.Range("G5").Value = "Hawaii"

# Nudge column E (the state column) to its new custom width
$ws.Columns("E").ColumnWidth = 19.5

# Move/save the active cell selection to G2
$ws.Range("G2").Select()
